$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value2
$newText = $oldText -replace [regex]::Escape("1000 Bs = 1.79 = 6513.43 pesos"), "1000 Bs = 1.82 = 6613.28 pesos"
$newText = $newText -replace [regex]::Escape("6513.43 pesos = 1.78 = 945.19 Bs"), "6613.28 pesos = 1.81 = 956.95 Bs"
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 549.5
$wsTasas.Range("O10").Value = 3634
$wsTasas.Range("N12").Value = 3648.9
$wsTasas.Range("O12").Value = 527.999
